$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "29.909.19"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  +1.45%  "

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "1.623.95"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  +1.19%  "

$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "0.992"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "  -0.64%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "213.60"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +0.47%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.519"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -0.47%  "

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.992"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  -0.70%  "

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "29.26"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  +9.02%  "

$ws.Range("E9").Value = "  +3.34%  "

$ws.Range("E10").Value = "  +0.89%  "

$ws.Range("E11").Value = "  +0.18%  "

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "1.854.11"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  +1.05%  "

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "1.616.12"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  +0.65%  "

$ws.Range("E14").Value = "  +6.19%  "

$ws.Range("E15").Value = "  +5.56%  "

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "29.916.10"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  +1.49%  "

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "8.83"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  +15.94%  "

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "64.39"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  +1.47%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "242.10"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +0.54%  "

$ws.Range("E20").Value = "  +2.34%  "

$ws.Range("E21").Value = "  -0.53%  "

$ws.Range("E22").Value = "  +2.76%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "9.60"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +4.30%  "

$ws.Range("E24").Value = "  +0.97%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "156.44"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +1.38%  "

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "15.62"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  +2.29%  "

$ws.Range("E27").Value = "  +1.13%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "6.59"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  +3.15%  "

$ws.Range("E29").Value = "  -0.64%  "

$ws.Range("E30").Value = "  +3.24%  "

$ws.Range("E31").Value = "  +5.13%  "

$ws.Range("E32").Value = "  +3.75%  "

$ws.Range("E33").Value = "  +4.29%  "

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "1.421.15"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  +0.71%  "

$ws.Range("E35").Value = "  +7.09%  "

$ws.Range("E36").Value = "  -0.22%  "

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "2.86"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  +1.03%  "

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "2.28"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  -0.78%  "

$ws.Range("E39").Value = "  +2.62%  "

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.556"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  +3.36%  "

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.0504"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +4.17%  "

$ws.Range("E42").Value = "  -0.17%  "

$ws.Range("E43").Value = "  +3.91%  "

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "69.26"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +5.04%  "

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "53.38"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  +0.22%  "

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "1.01"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  +18.42%  "

$ws.Range("E47").Value = "  -0.70%  "

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "5.43"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  +2.99%  "

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "1.763.40"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  +1.14%  "

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "88.47"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +2.10%  "

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.0₆0106"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +5.39%  "

